# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1117
$ws1.Range("F7").Value  = 274
$ws1.Range("F8").Value  = 49
$ws1.Range("F10").Value = 16129
$ws1.Range("F11").Value = 267
$ws1.Range("F14").Value = 6301
$ws1.Range("F21").Value = 12
$ws1.Range("F24").Value = 27
$ws1.Range("F35").Value = 140
$ws1.Range("F38").Value = 268

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1117
$ws4.Range("F7").Value  = 274
$ws4.Range("F8").Value  = 49
$ws4.Range("F10").Value = 16129
$ws4.Range("F11").Value = 267
$ws4.Range("F14").Value = 6301
$ws4.Range("F21").Value = 12
$ws4.Range("F24").Value = 27
$ws4.Range("F36").Value = 140
$ws4.Range("F39").Value = 268
